# QA and update files to console
# Applies the content/formatting edits described in the commit to the
# "ch" worksheet of the bindDialog workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 / Column C: tweak wording of the security-group description.
# Keep the red highlight + wrap formatting that the reviewer applied.
# ---------------------------------------------------------------------
$c2 = $ws.Range("C2")
$c2.Value = "Security group is a kind of distributed and statefull virtual firewall, which can be used by users to implement the network "
$c2.ClearFormats()
$c2.Font.Name = "Arial"
$c2.Font.Size = 12
$c2.Font.Color = 255
$c2.WrapText = $true

# ---------------------------------------------------------------------
# Row 4 / Column C: "Virtual Machines" -> "Virtual Machine" (singular),
# highlighted in red, no wrap.
# ---------------------------------------------------------------------
$c4 = $ws.Range("C4")
$c4.Value = "Virtual Machine"
$c4.ClearFormats()
$c4.Font.Name = "Arial"
$c4.Font.Size = 12
$c4.Font.Color = 255

# ---------------------------------------------------------------------
# Row 10 / Column C: lower-case the word "displayed", highlighted in red,
# no wrap.
# ---------------------------------------------------------------------
$c10 = $ws.Range("C10")
$c10.Value = "Configure the displayed tab"
$c10.ClearFormats()
$c10.Font.Name = "Arial"
$c10.Font.Size = 12
$c10.Font.Color = 255

# ---------------------------------------------------------------------
# Misc bookkeeping the reviewer left behind: the active selection moved
# to C15, and the page was switched to portrait orientation.
# ---------------------------------------------------------------------
[void]$ws.Range("C15").Select()
$ws.PageSetup.Orientation = 1
